$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- remove Sheet2 (merged into a single sheet workbook) ---
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

# --- widen column A to fit the new labels ---
$ws.Columns.Item(1).ColumnWidth = 25.33

# --- new AutoFill sample rows (matches the original authoring order) ---
$ws.Range("A8").Value = "Numbers"
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = "Time"
$ws.Range("B9").Value = 0.375

$ws.Range("A10").Value = "Week day"
$ws.Range("A12").Value = "Month"
$ws.Range("A14").Value = "Year"
$ws.Range("B14").Value = 2007

$ws.Range("B10").Value = "Monday"
$ws.Range("B13").Value = "Jan"
$ws.Range("A13").Value = "Month - short"
$ws.Range("A11").Value = "Week day - short"
$ws.Range("B12").Value = "Janurary"
$ws.Range("B11").Value = "Mon"

# --- row height + bold/18pt caption style for column A labels ---
$labels = $ws.Range("A8:A14")
$labels.RowHeight = 24
$labels.Font.Bold = $true
$labels.Font.Size = 18

# --- time-of-day number format for the Time sample ---
$ws.Range("B9").NumberFormat = "h:mm"

# --- selection matches the author's last edit position ---
$ws.Range("B12").Select()
